# Commit: "Added Zev's simulation with 100k voters"
# Observed change: sheet view zoom + selection updated, and explicit
# column widths set (narrow index column A, wider data columns B:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Custom column widths.
# Column A narrow (index/row-label column).
$ws.Columns.Item(1).ColumnWidth = 2.5
# Columns B:H wider (data columns) - set via EntireColumn on a multi-column range.
$ws.Range("B1:H1").EntireColumn.ColumnWidth = 6.666666666666667

# Update the selection to the full column H.
$ws.Range("H1:H1048576").Select() | Out-Null

# Zoom the sheet view to 162%.
$excel.ActiveWindow.Zoom = 162
